# Natmi Lif-Lifr output was regenerated following Dr Hou's advice: the
# ligand/receptor specificity values are recomputed (new normalisation),
# and the Sending/Target cluster grid is filled out completely (3 clusters
# x 3 clusters = 9 rows instead of the original 6 partial combinations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lif", "Lifr", "ECs", 1, 0.3333333333333333, 0.5355756666666667, 1.606727, 0.168705706571243, 0.168705706571243, 3, 1, 35.337883, 106.013649, 0.3968231145247413, 0.3968231145247413, 18.92611024631367, 170.334992216823, 0.06694632391969778, 0.06694632391969778),
    @("ECs", "Lif", "Lifr", "FAPs", 1, 0.3333333333333333, 0.5355756666666667, 1.606727, 0.168705706571243, 0.168705706571243, 3, 1, 40.49537033333333, 121.486111, 0.4547385869013757, 0.4547385869013756, 21.68833496318856, 195.195014668697, 0.07671699460840518, 0.07671699460840517),
    @("ECs", "Lif", "Lifr", "sCs", 1, 0.3333333333333333, 0.5355756666666667, 1.606727, 0.168705706571243, 0.168705706571243, 3, 1, 13.218724, 39.656172, 0.1484382985738831, 0.148438298573883, 7.079626918782667, 63.71664226904399, 0.02504238804314008, 0.02504238804314008),
    @("FAPs", "Lif", "Lifr", "ECs", 3, 1, 2.335066666666667, 7.0052, 0.7355432601013561, 0.7355432601013562, 3, 1, 35.337883, 106.013649, 0.3968231145247413, 0.3968231145247413, 82.51631266386666, 742.6468139748, 0.291880567341102, 0.291880567341102),
    @("FAPs", "Lif", "Lifr", "FAPs", 3, 1, 2.335066666666667, 7.0052, 0.7355432601013561, 0.7355432601013562, 3, 1, 40.49537033333333, 121.486111, 0.4547385869013757, 0.4547385869013756, 94.55938941968888, 851.0345047772, 0.3344799027033217, 0.3344799027033217),
    @("FAPs", "Lif", "Lifr", "sCs", 3, 1, 2.335066666666667, 7.0052, 0.7355432601013561, 0.7355432601013562, 3, 1, 13.218724, 39.656172, 0.1484382985738831, 0.148438298573883, 30.86660178826667, 277.7994160944, 0.1091827900569324, 0.1091827900569324),
    @("sCs", "Lif", "Lifr", "ECs", 3, 1, 0.3039726666666667, 0.911918, 0.09575103332740086, 0.09575103332740086, 3, 1, 35.337883, 106.013649, 0.3968231145247413, 0.3968231145247413, 10.74175052986467, 96.675754768782, 0.03799622326394152, 0.03799622326394151),
    @("sCs", "Lif", "Lifr", "FAPs", 3, 1, 0.3039726666666667, 0.911918, 0.09575103332740086, 0.09575103332740086, 3, 1, 40.49537033333333, 121.486111, 0.4547385869013757, 0.4547385869013756, 12.30948570787756, 110.785371370898, 0.0435416895896488, 0.04354168958964878),
    @("sCs", "Lif", "Lifr", "sCs", 3, 1, 0.3039726666666667, 0.911918, 0.09575103332740086, 0.09575103332740086, 3, 1, 13.218724, 39.656172, 0.1484382985738831, 0.148438298573883, 4.018130784210666, 36.163177057896, 0.01421312047381056, 0.01421312047381055)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}
